# Refresh cryptos list: updated Price/Volume(1h) figures, and corrected
# the swapped WEMIXToken/Dai rows (25 and 26) to their proper coin/link/values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.193.38"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "2.180.26"
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "66.54"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.34%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "36.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -12.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0930"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.58%  "
$ws.Range("D15").Value = "2.505.55"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.854"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "2.159.36"
$ws.Range("E18").Value = "  -3.01%  "
$ws.Range("D19").Value = "41.163.58"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("E20").Value = "  -2.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "228.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.69%  "
$ws.Range("E28").Value = "  -4.57%  "
$ws.Range("E29").Value = "  -3.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("E32").Value = "  -3.20%  "
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.41%  "
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.121"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.38%  "
$ws.Range("E37").Value = "  -4.41%  "
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.86%  "
$ws.Range("E40").Value = "  +3.24%  "
$ws.Range("E41").Value = "  -3.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.25%  "
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("E47").Value = "  -8.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("E51").Value = "  -3.77%  "
